# Glue_Board_BOM.xlsx update — "Added harnessing and assembly diagrams"
# The BOM rows that used to combine D2+D3 and SW1+SW2 into single lines are
# split into individual per-reference rows (each now carries its own MPN),
# which shifts every following row down. Two new rows are inserted to make
# room, the split rows are (re)populated, the Glue Board connector's MPN
# cell gets left-aligned, and the workbook-level bookkeeping (defined name
# range, dimension, selection) is refreshed to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Split the "D2 D3 " row (row 10) into separate "D2" / "D3 " rows ---
$ws.Rows("11").Insert()

$ws.Range("A10").Value = "D2"
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = "LED"
$ws.Range("D10").Value = "Connector_PinHeader_2.54mm:PinHeader_1x02_P2.54mm_Vertical"
$ws.Range("E10").Value = "~"
$ws.Range("F10").Value = "MPR3BWD"

$ws.Range("A11").Value = "D3 "
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = "LED"
$ws.Range("D11").Value = "Connector_PinHeader_2.54mm:PinHeader_1x02_P2.54mm_Vertical"
$ws.Range("E11").Value = "~"
$ws.Range("F11").Value = "MPR3RD"

# --- J6 (now row 17): left-align the Molex part-number cell ---
$ws.Range("F17").HorizontalAlignment = -4131

# --- Split the "SW1 SW2 " row (now row 35) into separate "SW1" / "SW2" rows ---
$ws.Rows("36").Insert()

$ws.Range("A35").Value = "SW1"
$ws.Range("B35").Value = 1
$ws.Range("C35").Value = "SW_SPST"
$ws.Range("D35").Value = "TerminalBlock:TerminalBlock_bornier-2_P5.08mm"
$ws.Range("E35").Value = "~"
$ws.Range("F35").Value = "M2011LL1W01-G"

$ws.Range("A36").Value = "SW2"
$ws.Range("B36").Value = 1
$ws.Range("C36").Value = "SW_SPST"
$ws.Range("D36").Value = "TerminalBlock:TerminalBlock_bornier-2_P5.08mm"
$ws.Range("E36").Value = "~"
$ws.Range("F36").Value = "M2011LL1W01-C"

# --- Defined name range grows from the 2 inserted rows ---
$wb.Names("Glue_Board").RefersToR1C1 = "=Sheet1!R1C1:R45C6"

# --- View bookkeeping: window geometry + active selection ---
$excel.Width = 1536
$excel.Height = 850
$ws.Cells.Item(1,1).Select()
$ws.Range("F37").Select()
